$d = $word.ActiveDocument

# --- Add the two new character styles referenced by the inserted code block ---
$defaultParaFont = $d.Styles.Item("DefaultParagraphFont")

$kw = $d.Styles.Add("hljs-keyword", 2)
$kw.BaseStyle = $defaultParaFont

$str = $d.Styles.Add("hljs-string", 2)
$str.BaseStyle = $defaultParaFont

# --- Locate the insertion point: end of the last real paragraph (the one that
#     holds the "_GoBack" bookmark), so its trailing bookmark markers stay at
#     the tail of the merged paragraph. ---
$target = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertRange = $d.Range($target.Range.End, $target.Range.End)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
        <w:body>
<w:p/><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="171717"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="171717"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>Use the </w:t></w:r><w:hyperlink r:id="rId7" w:history="1"><w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>Secret Manager tool</w:t></w:r></w:hyperlink><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="171717"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t> to set a password for these accounts.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="171717"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr></w:pPr></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rStyle w:val="hljs-keyword"/><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="0101FD"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="F2F2F2"/></w:rPr><w:t>dotnet</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="171717"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="F2F2F2"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rStyle w:val="hljs-keyword"/><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="0101FD"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="F2F2F2"/></w:rPr><w:t>user-secrets</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="171717"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="F2F2F2"/></w:rPr><w:t xml:space="preserve"> set </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="171717"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="F2F2F2"/></w:rPr><w:t>SeedUserPW</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="171717"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="F2F2F2"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rStyle w:val="hljs-string"/><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="A31515"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="F2F2F2"/></w:rPr><w:t>&lt;PW&gt;</w:t></w:r></w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
  <pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml">
    <pkg:xmlData>
      <Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">
        <Relationship Id="rId7" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://docs.microsoft.com/en-us/aspnet/core/security/app-secrets?view=aspnetcore-6.0&tabs=windows" TargetMode="External"/>
      </Relationships>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$insertRange.InsertXML($xml)

Write-Output "edit applied"
